$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: ID / TIME columns were leaking into the training features.
# Retrained + re-assessed the models; append the two new result rows
# (LogisticRegression, RandomForestClassifier) produced by the rerun.

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "LogisticRegression"
$ws.Range("C8").Value = 0.8711398394151816
$ws.Range("D8").Value = "{'clf__max_iter': 89}"
$ws.Range("E8").Value = "2023-09-26 10:44:12"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "RandomForestClassifier"
$ws.Range("C9").Value = 0.8873473226603699
$ws.Range("D9").Value = "{'clf__max_depth': 14, 'clf__n_estimators': 89}"
$ws.Range("E9").Value = "2023-09-26 10:45:00"

# Match the existing table's formatting: column A on data rows carries the
# bold/bordered/centered style (style index 1) used by A2:A7. Copy it over
# via PasteSpecial (formats only) so we reuse the existing style instead of
# minting a new one.
$ws.Range("A2").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
